# Delete the data row for "RGN" (Yangon, Myanmar), shifting all subsequent
# rows up by one. This is row 242 in the original sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(242).Delete()
